# "Ran a few tests" - bump the prior variances (column E, rows 2-22) from
# 0.05 to 0.15, switch the workbook to manual calculation (as Excel does
# while you're iterating on inputs without wanting a recalc each time),
# and leave the selection on E13 with the view scrolled back to the top.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Switch calculation mode to manual -> <calcPr calcMode="manual" .../>
$xlCalculationManual = -4135
$excel.Calculation = $xlCalculationManual

# Bump the prior variance column (E2:E22) from 0.05 to 0.15
$ws.Range("E2:E22").Value = 0.15

# Make sure this is the sheet we're looking at, scroll back to the top
# and leave the selection on E13
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("E13").Select()
